$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.608.65'
$ws.Range('E2').Value = '  +0.88%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.803.74'
$ws.Range('E3').Value = '  -0.39%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.56'
$ws.Range('E5').Value = '  -0.28%  '

$ws.Range('E6').Value = '  -0.05%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5452'
$ws.Range('E7').Value = '  -4.40%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3794'
$ws.Range('E8').Value = '  -2.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07514'
$ws.Range('E9').Value = '  -1.22%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.38'
$ws.Range('E10').Value = '  -1.76%  '

$ws.Range('E11').Value = '  -2.24%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.000'
$ws.Range('E12').Value = '  -0.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.66'
$ws.Range('E13').Value = '  -2.57%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.154'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.382'
$ws.Range('E15').Value = '  +1.62%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.803.00'
$ws.Range('E16').Value = '  -0.23%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '90.16'
$ws.Range('E17').Value = '  -1.97%  '

$ws.Range('E18').Value = '  -0.77%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06480'
$ws.Range('E19').Value = '  -0.01%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.09%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.31'
$ws.Range('E21').Value = '  -0.06%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.936'
$ws.Range('E22').Value = '  -1.21%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.641.27'
$ws.Range('E23').Value = '  +0.96%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.11'
$ws.Range('E24').Value = '  -1.74%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.098'
$ws.Range('E25').Value = '  -2.02%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.58'
$ws.Range('E26').Value = '  +1.69%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.44'
$ws.Range('E27').Value = '  -1.85%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.010.30'
$ws.Range('E28').Value = '  -0.39%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.358'
$ws.Range('E29').Value = '  -3.59%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '123.02'
$ws.Range('E30').Value = '  -0.85%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.108'
$ws.Range('E31').Value = '  -4.45%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1056'
$ws.Range('E32').Value = '  -0.85%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.632'
$ws.Range('E33').Value = '  -2.65%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.681'
$ws.Range('E34').Value = '  +1.42%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.06625'
$ws.Range('E35').Value = '  +8.31%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2253'
$ws.Range('E36').Value = '  +2.16%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02302'
$ws.Range('E37').Value = '  -0.83%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.755'
$ws.Range('E38').Value = '  -1.81%  '

$ws.Range('E39').Value = '  -0.56%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6248'
$ws.Range('E40').Value = '  -2.51%  '

$ws.Range('E41').Value = '  -3.71%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.196'
$ws.Range('E42').Value = '  +2.84%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.436'
$ws.Range('E43').Value = '  +4.16%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9995'
$ws.Range('E44').Value = '  -0.11%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.24'
$ws.Range('E45').Value = '  -1.55%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.697'
$ws.Range('E46').Value = '  -0.22%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5850'
$ws.Range('E47').Value = '  -2.61%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.67'
$ws.Range('E48').Value = '  +3.81%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.945'
$ws.Range('E49').Value = '  -0.02%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.160'
$ws.Range('E50').Value = '  +0.93%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06886'
$ws.Range('E51').Value = '  +0.27%  '
